$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that often LOOKS like a plain number
# (e.g. "106.00", "3.70", "1.937.40"). A bare .Value assignment lets Excel
# auto-coerce numeric-looking text into a real Number, which loses trailing
# zeros / exact text and flips the cell type from Text to Number. Forcing
# the cell to Text format ("@") before the write keeps it a literal string;
# resetting the style back to "Normal" afterwards avoids leaving a stray
# number-format style behind (the source workbook never assigns a style to
# these data cells).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '34.144.62'
$ws.Range("E2").Value = '  -1.35%  '
Set-TextValue $ws.Range("D3") '1.783.47'
$ws.Range("E3").Value = '  -0.95%  '
Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  +0.13%  '
Set-TextValue $ws.Range("D5") '221.41'
$ws.Range("E5").Value = '  -2.05%  '
Set-TextValue $ws.Range("D6") '0.549'
$ws.Range("E7").Value = '  +0.18%  '
Set-TextValue $ws.Range("D8") '31.67'
$ws.Range("E8").Value = '  -4.08%  '
Set-TextValue $ws.Range("D9") '0.288'
$ws.Range("E9").Value = '  +1.16%  '
Set-TextValue $ws.Range("D10") '0.0709'
$ws.Range("E10").Value = '  +6.37%  '
Set-TextValue $ws.Range("D11") '0.0924'
$ws.Range("E11").Value = '  -0.92%  '
Set-TextValue $ws.Range("D12") '2.040.93'
$ws.Range("E12").Value = '  -0.84%  '
Set-TextValue $ws.Range("D13") '1.787.91'
$ws.Range("E13").Value = '  -0.52%  '
Set-TextValue $ws.Range("D14") '10.66'
$ws.Range("E14").Value = '  -4.12%  '
Set-TextValue $ws.Range("D15") '0.625'
Set-TextValue $ws.Range("D16") '34.107.13'
$ws.Range("E16").Value = '  -1.31%  '
Set-TextValue $ws.Range("D17") '4.22'
$ws.Range("E17").Value = '  -1.44%  '
Set-TextValue $ws.Range("D18") '67.88'
$ws.Range("E18").Value = '  -2.70%  '
Set-TextValue $ws.Range("D19") '244.91'
$ws.Range("E19").Value = '  -4.39%  '
Set-TextValue $ws.Range("D20") '0.0₃0778'
$ws.Range("E20").Value = '  +2.53%  '
$ws.Range("E21").Value = '  +0.22%  '
Set-TextValue $ws.Range("D22") '10.65'
$ws.Range("E22").Value = '  +1.47%  '
Set-TextValue $ws.Range("D23") '4.10'
$ws.Range("E23").Value = '  -3.45%  '
Set-TextValue $ws.Range("D24") '2.12'
$ws.Range("E24").Value = '  -0.35%  '
Set-TextValue $ws.Range("D25") '157.66'
$ws.Range("E25").Value = '  -0.16%  '
Set-TextValue $ws.Range("D26") '16.36'
$ws.Range("E26").Value = '  -1.00%  '
Set-TextValue $ws.Range("D27") '7.04'
$ws.Range("E28").Value = '  -1.95%  '
Set-TextValue $ws.Range("D29") '0.999'
$ws.Range("E29").Value = '  +0.12%  '
Set-TextValue $ws.Range("D30") '0.0519'
$ws.Range("E30").Value = '  +0.27%  '
Set-TextValue $ws.Range("D31") '3.70'
$ws.Range("E31").Value = '  -2.63%  '
Set-TextValue $ws.Range("D32") '1.20'
$ws.Range("E32").Value = '  +0.66%  '
Set-TextValue $ws.Range("D33") '3.52'
$ws.Range("E33").Value = '  -2.44%  '
Set-TextValue $ws.Range("D34") '1.84'
$ws.Range("E34").Value = '  -4.12%  '
Set-TextValue $ws.Range("D35") '1.395.44'
$ws.Range("E35").Value = '  -4.61%  '
$ws.Range("E36").Value = '  -0.44%  '
Set-TextValue $ws.Range("D37") '0.628'
$ws.Range("E37").Value = '  -0.97%  '
Set-TextValue $ws.Range("D38") '0.0186'
$ws.Range("E38").Value = '  -2.11%  '
Set-TextValue $ws.Range("D39") '2.77'
$ws.Range("E39").Value = '  -2.99%  '
Set-TextValue $ws.Range("D40") '0.933'
$ws.Range("E40").Value = '  +3.50%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D41") '79.55'
$ws.Range("E41").Value = '  -4.53%  '
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range("D42") '2.35'
$ws.Range("E42").Value = '  +1.34%  '
$ws.Range("E43").Value = '  +0.87%  '
Set-TextValue $ws.Range("D44") '0.0495'
$ws.Range("E44").Value = '  -2.49%  '
Set-TextValue $ws.Range("D45") '1.04'
$ws.Range("E45").Value = '  +0.37%  '
Set-TextValue $ws.Range("D46") '5.84'
$ws.Range("E46").Value = '  -1.27%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range("D47") '106.00'
$ws.Range("E47").Value = '  +5.26%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range("D48") '1.937.40'
$ws.Range("E48").Value = '  -0.96%  '
Set-TextValue $ws.Range("D49") '0.995'
$ws.Range("E49").Value = '  -0.32%  '
Set-TextValue $ws.Range("D50") '11.80'
$ws.Range("E50").Value = '  -1.27%  '
Set-TextValue $ws.Range("D51") '0.0₆0120'
$ws.Range("E51").Value = '  +3.32%  '
